$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting (styles) from matching existing rows ---
$ws.Range("A3:E3").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("A10:E10").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A40:E40").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A41:E41").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A42:E42").PasteSpecial(-4122)

# --- Set cell values ---
$ws.Range("A24").Value = 'SCRIPT/P01P04A/um2503.ssb'

$ws.Range("A25").Value = 'SCRIPT/P01P04A/us0103.ssb'
$ws.Range("B25").Value = 196
$ws.Range("C25").Value = ' It\''s thanks to you guys that I\''m\neven able to sit here enjoying\nsuch delicious drinks.'
$ws.Range("D25").Value = ' Друзья, благодаря вам я могу\nсидеть здесь и наслаждаться вкуснейшими\nнапитками.'
$ws.Range("E25").Value = ' Äñôèûÿ, áìàãïäàñÿ âàí ÿ íïãô\nòéäåóû èäåòû é îàòìàçäàóûòÿ âëôòîåêšéíé\nîàðéóëàíé.'

$ws.Range("B26").Value = 199
$ws.Range("C26").Value = ' Thanks, Team [team:].'
$ws.Range("D26").Value = ' Спасибо вам, Команда\n[team:].'
$ws.Range("E26").Value = ' Òðàòéáï âàí, Ëïíàîäà\n[team:].'

$ws.Range("A27").Value = 'SCRIPT/P01P04A/us3108.ssb'
$ws.Range("B27").Value = 177
$ws.Range("C27").Value = ' I wonder if any new drink\ningredients can be found around [CS:P]Sky Peak[CR]…'
$ws.Range("D27").Value = ' Интересно, можно ли на [CS:P]Небесном\nПике[CR] найти новых ингредиентов для\nнапитков...'
$ws.Range("E27").Value = ' Éîóåñåòîï, íïçîï ìé îà [CS:P]Îåáåòîïí\nÐéëå[CR] îàêóé îïâúö éîãñåäéåîóïâ äìÿ\nîàðéóëïâ…'

$ws.Range("A28").Value = 'SCRIPT/D73P23A/us3107.ssb'
$ws.Range("B28").Value = 112
$ws.Range("C28").Value = ' I heard [CS:N]Shaymin[CR] talking earlier…'
$ws.Range("D28").Value = ' Я подслушал, о чём говорил\n[CS:N]Шеймин[CR]...'
$ws.Range("E28").Value = ' Ÿ ðïäòìôšàì, ï œæí ãïâïñéì\n[CS:N]Šåêíéî[CR]...'

$ws.Range("B29").Value = 115
$ws.Range("C29").Value = ' I did find a [CS:I]Sky Gift[CR]…'
$ws.Range("D29").Value = ' Я нашёл [CS:I]Небесный Подарок[CR]...'
$ws.Range("E29").Value = ' Ÿ îàšæì [CS:I]Îåáåòîúê Ðïäàñïë[CR]...'

$ws.Range("B30").Value = 118
$ws.Range("C30").Value = ' Team [team:] saved the\nworld, and I owe you for many things, so…'
$ws.Range("D30").Value = ' Команда [team:] спасла\nмир и я вам стольким обязан...'
$ws.Range("E30").Value = ' Ëïíàîäà [team:] òðàòìà\níéñ é ÿ âàí òóïìûëéí ïáÿèàî...'

$ws.Range("B31").Value = 121
$ws.Range("C31").Value = ' I always wanted to show you\nmy appreciation.'
$ws.Range("D31").Value = ' Я всегда хотел выразить свою\nпризнательность.'
$ws.Range("E31").Value = ' Ÿ âòåãäà öïóåì âúñàèéóû òâïý\nðñéèîàóåìûîïòóû.'

$ws.Range("B32").Value = 124
$ws.Range("C32").Value = ' So please accept this from me!'
$ws.Range("D32").Value = ' Поэтому, прошу, прими его в\nдар!'
$ws.Range("E32").Value = ' Ðïüóïíô, ðñïšô, ðñéíé åãï â\näàñ!'

$ws.Range("B33").Value = 134
$ws.Range("C33").Value = '[CN][player] received\n[CN]a [CS:I]Sky Gift[CR].'
$ws.Range("D33").Value = '[CN][player] получает\n[CN][CS:I]Небесный Подарок[CR].'
$ws.Range("E33").Value = '[CN][player] ðïìôœàåó\n[CN][CS:I]Îåáåòîúê Ðïäàñïë[CR].'

$ws.Range("B34").Value = 137
$ws.Range("C34").Value = '[CN]You gently open the lid…'
$ws.Range("D34").Value = '[CN]Вы осторожно открываете его...'
$ws.Range("E34").Value = '[CN]Âú ïòóïñïçîï ïóëñúâàåóå åãï...'

$ws.Range("B35").Value = 143
$ws.Range("C35").Value = '[CN]Inside was the\n[CN][s_item:0]!'
$ws.Range("D35").Value = '[CN]Внутри лежит предмет\n[CN][s_item:0]!'
$ws.Range("E35").Value = '[CN]Âîôóñé ìåçéó ðñåäíåó\n[CN][s_item:0]!'

$ws.Range("B36").Value = 146
$ws.Range("C36").Value = '[CN]A warm and fuzzy feeling falls over you…'
$ws.Range("D36").Value = '[CN]Вас охватывает тёплое и приятное чувство...'
$ws.Range("E36").Value = '[CN]Âàò ïöâàóúâàåó óæðìïå é ðñéÿóîïå œôâòóâï...'

$ws.Range("B37").Value = 153
$ws.Range("C37").Value = ' ...Or so I thought, but you seem\nto have too many items already.'
$ws.Range("D37").Value = ' ...Или я хотел, но, похоже, у\nтебя слишком много вещей.'
$ws.Range("E37").Value = ' ...Éìé ÿ öïóåì, îï, ðïöïçå, ô\nóåáÿ òìéšëïí íîïãï âåþåê.'

$ws.Range("B38").Value = 156
$ws.Range("C38").Value = ' Hmm... That\''s too bad…'
$ws.Range("D38").Value = ' Хмм... Очень жаль...'
$ws.Range("E38").Value = ' Öíí... Ïœåîû çàìû...'

$ws.Range("B39").Value = 96
$ws.Range("C39").Value = ' Huh? Thanks for the gift?'
$ws.Range("D39").Value = ' Что? Спасибо за подарок?'
$ws.Range("E39").Value = ' Œóï? Òðàòéáï èà ðïäàñïë?'

$ws.Range("B40").Value = 99
$ws.Range("C40").Value = ' You\''re more than welcome!'
$ws.Range("D40").Value = ' Всегда пожалуйста!'
$ws.Range("E40").Value = ' Âòåãäà ðïçàìôêòóà!'

$ws.Range("A41").Value = 'SCRIPT/D73P24A/us3104.ssb'
$ws.Range("B41").Value = 79
$ws.Range("C41").Value = ' Phew... The 4th Station\nClearing, huh?'
$ws.Range("D41").Value = ' Фух... Поляна 4-го Перехода, да?'
$ws.Range("E41").Value = ' Õôö... Ðïìÿîà 4-ãï Ðåñåöïäà, äà?'

$ws.Range("B42").Value = 82
$ws.Range("C42").Value = ' As I find a [CS:I]Sky Gift[CR], I search\nfor someone to give it to... I\''m not\nmaking much progress.'
$ws.Range("D42").Value = ' Как только я нахожу [CS:I]Небесный\nПодарок[CR], я пытаюсь найти кого-нибудь, кому\nмогу его дать... У меня не получается.'
$ws.Range("E42").Value = ' Ëàë óïìûëï ÿ îàöïçô [CS:I]Îåáåòîúê\nÐïäàñïë[CR], ÿ ðúóàýòû îàêóé ëïãï-îéáôäû, ëïíô\níïãô åãï äàóû... Ô íåîÿ îå ðïìôœàåóòÿ.'

# --- Set row heights ---
$ws.Rows.Item(24).RowHeight = 43.2
$ws.Rows.Item(25).RowHeight = 43.2
$ws.Rows.Item(27).RowHeight = 43.2
$ws.Rows.Item(28).RowHeight = 43.2
$ws.Rows.Item(30).RowHeight = 21.6
$ws.Rows.Item(31).RowHeight = 21.6
$ws.Rows.Item(33).RowHeight = 31.8
$ws.Rows.Item(35).RowHeight = 21.6
$ws.Rows.Item(36).RowHeight = 21.6
$ws.Rows.Item(37).RowHeight = 21.6
$ws.Rows.Item(41).RowHeight = 31.2
$ws.Rows.Item(42).RowHeight = 42

# --- Update sheet view ---
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("C44").Select()
